$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Snapshot the "before" values for columns A, B, F, G across rows 5-41 ---
# (F only exists on some rows; default to empty string when absent)
$rows = 5..41
$colA = @{}
$colB = @{}
$colF = @{}
$colG = @{}
foreach ($r in $rows) {
    $colA[$r] = $ws.Cells.Item($r, 1).Value2()
    $colB[$r] = $ws.Cells.Item($r, 2).Value2()
    $colF[$r] = $ws.Cells.Item($r, 6).Value2()
    $colG[$r] = $ws.Cells.Item($r, 7).Value2()
}

# --- Mapping: destination row -> source row (which row's A/B/F/G values it should receive) ---
$srcForDest = @{
    5 = 33
    6 = 17
    7 = 36
    8 = 35
    9 = 29
    10 = 6
    11 = 32
    12 = 20
    13 = 15
    14 = 16
    15 = 14
    16 = 25
    17 = 5
    18 = 18
    19 = 19
    20 = 8
    21 = 10
    22 = 9
    23 = 11
    24 = 12
    25 = 31
    26 = 21
    27 = 22
    28 = 26
    29 = 13
    30 = 27
    31 = 37
    32 = 34
    33 = 38
    34 = 24
    35 = 7
    36 = 23
    37 = 30
    38 = 28
    39 = 41
    40 = 40
    41 = 39
}

foreach ($dest in $srcForDest.Keys) {
    $src = $srcForDest[$dest]
    $ws.Cells.Item($dest, 1).Value = $colA[$src]
    $ws.Cells.Item($dest, 2).Value = $colB[$src]
    $ws.Cells.Item($dest, 7).Value = $colG[$src]
    $fval = $colF[$src]
    if ($fval) {
        $ws.Cells.Item($dest, 6).Value = $fval
    } else {
        $ws.Cells.Item($dest, 6).Value = ""
    }
}

# --- Update column C ("Förändrad") from 46062 to 46063 for all data rows (2-41) ---
foreach ($r in 2..41) {
    $ws.Cells.Item($r, 3).Value = 46063
}

"Done."
